# Investor workbook: insert a "Primary Email" column (with mailto hyperlinks)
# between the existing "PAN *" and "TAGS" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the width (in "characters") of column B so the freshly inserted
# column C can inherit the same look (cols A:C all end up width 35).
$sourceWidth = $ws.Columns("B:B").ColumnWidth

# Insert a new blank column at C; everything from old C onward shifts right.
$ws.Columns("C:C").Insert()
$ws.Columns("C:C").ColumnWidth = $sourceWidth

# Header for the new column.
$ws.Range("C1").Value = "Primary Email"

# Populate the new column with each investor's primary e-mail address and
# turn each one into a mailto: hyperlink (Excel auto-applies the built-in
# "Hyperlink" style to the cell).
$emails = @("emp1@gmail.com", "emp2@gmail.com", "emp3@gmail.com", "emp4@gmail.com", "emp5@gmail.com", "emp6@gmail.com")

for ($i = 0; $i -lt $emails.Count; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $emails[$i]
    $null = $ws.Hyperlinks.Add($cell, "mailto:" + $emails[$i])
}

$null = $ws.Range("C8").Select()
